$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (matches ExpiryDate 2025-08-29 / B12 anchor date) - fill in previously blank M1_PH/CM2_PH/CMN3_PH/CMN4_PH values
$ws.Range("I2").Value = -0.7766300884562578
$ws.Range("J2").Value = 0.197766400003976
$ws.Range("K2").Value = -0.2003848170313115
$ws.Range("L2").Value = 2.53868442540057

# Row 18 (matches ExpiryDate 2025-09-26) - fill in previously blank M1_PH/CM2_PH/CMN3_PH/CMN4_PH values
$ws.Range("I18").Value = -0.9902044561778961
$ws.Range("J18").Value = 0.2497106689931531
$ws.Range("K18").Value = 0.2162614919740294
$ws.Range("L18").Value = 2.124678745864521
